# ESP32-S3 WROOM 2 hinzugefügt.
# Add a new "Mouser" section to the BOM sheet: a sub-header in column D
# (mirroring the "Digikey" header in D8) followed by a new component row
# (quantity, component name, Mouser order number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New supplier sub-header (row 25, leaving row 24 blank as a separator).
$ws.Range("D25").Value = "Mouser"

# New BOM line (row 26): quantity, component, then order number so the
# shared-string table is populated in the same order as the source edit.
$ws.Range("A26").Value = 1
$ws.Range("D26").Value = "710-830003000 "
$ws.Range("B26").Value = "Crystal 32.786kHz"

# Selection ends up resting on the next empty row, same as after typing
# the new row in the Excel UI and pressing Enter/Down.
$ws.Range("A27").Select()
